# Apply the row-level data update described by the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data layout: employee_id, employee_name, department, absence_reason,
# absence_duration, absence_date, salary
$rows = @(
    @{ Row = 2;  A = 94947; B = "Maria Liz Nogueira";      C = "Recursos Humanos"; D = "Doenca";              E = 1; F = 45089; G = 7010.61 }
    @{ Row = 3;  A = 32206; B = "Pedro Henrique Vargas";   C = "Marketing";        D = "Consulta medica";     E = 8; F = 45103; G = 4000.13 }
    @{ Row = 4;  A = 81541; B = "Maria Fernanda Lopes";    C = "Operacoes";        D = "Outros";               E = 2; F = 45084; G = 3449.96 }
    @{ Row = 5;  A = 60115; B = "Thiago Fonseca";          C = "TI";               D = "Outros";               E = 5; F = 45079; G = 2288.44 }
    @{ Row = 6;  A = 25152; B = "Julia Dias";              C = "P&D";              D = "Consulta medica";     E = 8; F = 45084; G = 3214.49 }
    @{ Row = 7;  A = 23359; B = "Josué Cassiano";          C = "P&D";              D = "Viagem de negocios";  E = 8; F = 45096; G = 7037.41 }
    @{ Row = 8;  A = 80877; B = "Maria Cecília Azevedo";   C = "Engenharia";       D = "Consulta medica";     E = 7; F = 45104; G = 5106.57 }
    @{ Row = 9;  A = 92030; B = "Apollo Peixoto";          C = "TI";               D = "Problemas pessoais";  E = 3; F = 45099; G = 7201.75 }
    @{ Row = 10; A = 37909; B = "Carlos Eduardo Pastor";   C = "Financeiro";       D = "Viagem de negocios";  E = 6; F = 45088; G = 2907.44 }
    @{ Row = 11; A = 78802; B = "Antony Monteiro";         C = "Financeiro";       D = "Viagem de negocios";  E = 6; F = 45102; G = 6388.7 }
)

foreach ($r in $rows) {
    $rowNum = $r.Row
    $ws.Cells.Item($rowNum, 1).Value = $r.A
    $ws.Cells.Item($rowNum, 2).Value = $r.B
    $ws.Cells.Item($rowNum, 3).Value = $r.C
    $ws.Cells.Item($rowNum, 4).Value = $r.D
    $ws.Cells.Item($rowNum, 5).Value = $r.E
    $ws.Cells.Item($rowNum, 6).Value = $r.F
    $ws.Cells.Item($rowNum, 7).Value = $r.G
}
